# edit.ps1 - apply the "changed (a little bit) personal description" edit
#
# Summary of changes (per the unified OOXML diff):
#  1. ". I live in Petrozavodsk, Karelia Republic. I am on 11th grade in
#     school now." -> split into three runs ending with
#     ". I live in Petrozavodsk, Karelia Republic. I finished 11 grades in
#     school."
#  2. "I've" -> split into "I'" + a (now-empty) _GoBack bookmark + "ve"
#  3. "Languages" + ":" (two runs) merged into a single run "Languages:"
#  4. The _GoBack bookmark at the very end of the document (after
#     "Learn new things") is removed.
#
# NOTE: this PowerShell-interop runtime does not bind named (-Param value)
# function arguments correctly, so helper functions below use plain
# positional parameters only.

function New-OpenXmlRun($Text, $PreserveSpace) {
    $space = ''
    if ($PreserveSpace) { $space = ' xml:space="preserve"' }
    $rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr>'
    return '<w:r>' + $rPr + '<w:t' + $space + '>' + $Text + '</w:t></w:r>'
}

function New-OpenXmlPackage($BodyXml) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $footer = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $BodyXml + $footer
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) ". I live in Petrozavodsk, Karelia Republic. I am on 11th grade in
#    school now." -> 3 runs
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(". I live in Petrozavodsk, Karelia Republic. I am on 11th grade in school now.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Petrozavodsk/grade sentence to replace"
}
$target = $d.Range($r.Start, $r.End)

$body = New-OpenXmlRun ". I live in Petrozavodsk, Karelia Republic. I " $true
$body += New-OpenXmlRun "finished 11 grades in school" $false
$body += New-OpenXmlRun "." $false
$target.InsertXML((New-OpenXmlPackage $body))

# ---------------------------------------------------------------------
# 2) "I've" (the first occurrence, right after the gramStart proofErr
#    before "started to learn programming") -> "I'" + _GoBack bookmark + "ve"
#    Only touch the "ve" portion so the "I'" remainder keeps being part
#    of the original (untouched) run / proofErr wrapping.
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("I've started to learn programming", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'I've started to learn programming'"
}
$iveStart = $r.Start
$veRange = $d.Range($iveStart + 2, $iveStart + 4)
if ($veRange.Text -ne "ve") {
    throw "Unexpected text at the I've split point: [$($veRange.Text)]"
}

$bodyVe = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + (New-OpenXmlRun "ve" $false)
$veRange.InsertXML((New-OpenXmlPackage $bodyVe))

# ---------------------------------------------------------------------
# 3) "Languages" + ":" (two runs) -> single run "Languages:"
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Languages:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Languages:'"
}
$target = $d.Range($r.Start, $r.End)
$bodyLang = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Languages:</w:t></w:r>'
$target.InsertXML((New-OpenXmlPackage $bodyLang))

# ---------------------------------------------------------------------
# 4) Remove the _GoBack bookmark at the end of the document (after
#    "Learn new things")
# ---------------------------------------------------------------------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
    # already gone / nothing to remove
}

Write-Output "done"
